$d = $word.ActiveDocument

# Remove the "nome_" prefix run from the "{{nome_responsavel_pela_demanda}}" placeholder,
# merging it away so the field reads "{{responsavel_pela_demanda}}"
$d.Content.Find.Execute("nome_", $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
